# Normalize column F (the timestamp column) for rows 4 through 108:
# every one of these cells gets set to the same datetime serial value
# that rows 2 and 3 already hold (45695.58159722222).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4:F108").Value = 45695.58159722222
